# Update the "dSF" column (F) values for a set of rows, as the data was
# repulled / recalculated (repull data, push all data, mean calculation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = 3
    11 = 0
    12 = 0
    14 = 3
    15 = 2
    33 = -3
    41 = 1
    43 = 2
    45 = 3
    46 = 0
    60 = -3
    64 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
